# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) timestamps on the
# per-language report sheets, as produced by a fresh handback report run.

$wb = $excel.ActiveWorkbook

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("E2").Value = "2016-03-21 20:21:06"
$wsZh.Range("E5").Value = "2016-03-21 20:21:06"
$wsZh.Range("H2").Value = "2016-03-21 20:21:43"
$wsZh.Range("H5").Value = "2016-03-21 20:21:43"

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("E2").Value = "2016-03-21 20:21:17"
$wsDe.Range("E5").Value = "2016-03-21 20:21:17"
$wsDe.Range("H2").Value = "2016-03-21 20:21:50"
$wsDe.Range("H5").Value = "2016-03-21 20:21:50"
